$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("i18n-de-progress")

# --- Data updates -----------------------------------------------------
# Rows 13-40: column E becomes "yes" (losing its date-format style),
# column F becomes "in progress" (keeping its existing style).
for ($r = 13; $r -le 40; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $eCell.ClearFormats()
    $eCell.Value2 = "yes"

    $fCell = $ws.Cells.Item($r, 6)   # column F
    $fCell.Value2 = "in progress"
}

# Rows 41-52: column F becomes "yes" (losing its date-format style).
# Column E in this range is already "yes" with no style and stays as-is.
for ($r = 41; $r -le 52; $r++) {
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $fCell.ClearFormats()
    $fCell.Value2 = "yes"
}

# --- View / selection state --------------------------------------------
$aw = $excel.ActiveWindow
$aw.FreezePanes = $true
[void]$ws.Range("A2").Select()
[void]$ws.Range("E41").Select()

$wb.Application.Calculate()
